$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-09-21 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-22 Sunday", 2)

# Update the division problems in the table. Using direct cell addressing
# (row, column) avoids any ambiguity from values that coincide with other
# cells' old/new text (e.g. "92÷2=" becomes "19÷5=" while a different cell's
# "19÷5=" becomes "33÷8=").
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "43÷2=" },
    @{ Row = 1;  Col = 2; Text = "47÷8=" },
    @{ Row = 1;  Col = 3; Text = "45÷2=" },
    @{ Row = 1;  Col = 4; Text = "29÷5=" },
    @{ Row = 1;  Col = 5; Text = "46÷2=" },

    @{ Row = 5;  Col = 1; Text = "32÷2=" },
    @{ Row = 5;  Col = 2; Text = "70÷2=" },
    @{ Row = 5;  Col = 3; Text = "57÷9=" },
    @{ Row = 5;  Col = 4; Text = "19÷5=" },
    @{ Row = 5;  Col = 5; Text = "97÷9=" },

    @{ Row = 9;  Col = 1; Text = "72÷2=" },
    @{ Row = 9;  Col = 2; Text = "64÷8=" },
    @{ Row = 9;  Col = 3; Text = "53÷5=" },
    @{ Row = 9;  Col = 4; Text = "69÷3=" },
    @{ Row = 9;  Col = 5; Text = "23÷7=" },

    @{ Row = 13; Col = 1; Text = "32÷4=" },
    @{ Row = 13; Col = 2; Text = "11÷7=" },
    @{ Row = 13; Col = 3; Text = "97÷4=" },
    @{ Row = 13; Col = 4; Text = "94÷5=" },
    @{ Row = 13; Col = 5; Text = "39÷9=" },

    @{ Row = 17; Col = 1; Text = "92÷5=" },
    @{ Row = 17; Col = 2; Text = "52÷3=" },
    @{ Row = 17; Col = 3; Text = "16÷5=" },
    @{ Row = 17; Col = 4; Text = "55÷3=" },
    @{ Row = 17; Col = 5; Text = "33÷8=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
